$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for the new "Start Page / Card Viewer" fragment row -------
# Insert 4 blank rows before row 4; this pushes the "Deck Manager..." flow
# (previously rows 6-14, plus the "Represents multi-pane grouping" note on
# row 4) down by four rows.
$ws.Rows("4:7").Insert()

# The old up/down "arrow" placeholder text in the header block is no longer
# needed in its old spots - clear it out.
$ws.Range("B3").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("B9").ClearContents()

# New header block (rows 3-6, bold like the "User" label on row 2) with the
# new fragment names.
$ws.Range("B3:B6").Font.Bold = $true
$ws.Range("B6").Value = "Start Page (dashboard)"
$ws.Range("D6").Value = "Card Viewer"
$ws.Range("D6").Font.Bold = $true

# Column widths nudge slightly now that the bold fragment labels are the
# widest entries in their columns.
$ws.Columns("A").ColumnWidth = 18.75
$ws.Columns("B").ColumnWidth = 22.5836
$ws.Columns("D").ColumnWidth = 17.25
$ws.Columns("G").ColumnWidth = 56.25

# --- Flow-diagram connector arrows between the new boxes ------------------
$c1 = $ws.Shapes.AddConnector(1, 222.75, 87.75, 82.5, 39.75)
$c1.Name = "Straight Arrow Connector 2"
$c1.Line.BeginArrowheadStyle = 2
$c1.Line.EndArrowheadStyle = 2

$c2 = $ws.Shapes.AddConnector(1, 168, 93, 0, 38.25)
$c2.Name = "Straight Arrow Connector 4"
$c2.Line.BeginArrowheadStyle = 2
$c2.Line.EndArrowheadStyle = 2

$c3 = $ws.Shapes.AddConnector(1, 168, 33.75, 0, 41.25)
$c3.Name = "Straight Arrow Connector 6"
$c3.Line.BeginArrowheadStyle = 2
$c3.Line.EndArrowheadStyle = 2

$c4 = $ws.Shapes.AddConnector(1, 330, 91.5, 0, 102)
$c4.Name = "Straight Arrow Connector 8"
$c4.Line.BeginArrowheadStyle = 2
$c4.Line.EndArrowheadStyle = 2

# Leave the selection where the author's last edit was.
$ws.Range("D6").Select()
